$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.410.19"
$ws.Range("E2").Value = "  -4.60%  "
$ws.Range("D3").Value = "1.568.31"
$ws.Range("E3").Value = "  -4.91%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "291.37"
$ws.Range("E6").Value = "  -2.79%  "
$ws.Range("D7").Value = "0.3691"
$ws.Range("E7").Value = "  -2.56%  "
$ws.Range("D8").Value = "49.64"
$ws.Range("E8").Value = "  -1.54%  "
$ws.Range("D9").Value = "0.3368"
$ws.Range("E9").Value = "  -5.45%  "
$ws.Range("D10").Value = "1.167"
$ws.Range("E10").Value = "  -4.44%  "
$ws.Range("D11").Value = "0.07559"
$ws.Range("E11").Value = "  -6.67%  "
$ws.Range("E12").Value = "  -0.08%  "
$ws.Range("D13").Value = "21.08"
$ws.Range("E13").Value = "  -4.37%  "
$ws.Range("D14").Value = "6.041"
$ws.Range("E14").Value = "  -5.69%  "
$ws.Range("D15").Value = "6.846"
$ws.Range("E15").Value = "  -7.32%  "
$ws.Range("D16").Value = "0.00001142"
$ws.Range("E16").Value = "  -4.45%  "
$ws.Range("D17").Value = "1.575.54"
$ws.Range("E17").Value = "  -5.00%  "
$ws.Range("D18").Value = "89.14"
$ws.Range("E18").Value = "  -8.27%  "
$ws.Range("D19").Value = "0.06705"
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "6.245"
$ws.Range("E21").Value = "  -7.38%  "
$ws.Range("B22").Value = "BitDAO"
$ws.Range("C22").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D22").Value = "0.5289"
$ws.Range("E22").Value = "  -8.66%  "
$ws.Range("B23").Value = "Avalanche"
$ws.Range("C23").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D23").Value = "16.34"
$ws.Range("E23").Value = "  -5.55%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "11.93"
$ws.Range("E24").Value = "  -4.03%  "
$ws.Range("B25").Value = "WrappedBTC"
$ws.Range("C25").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D25").Value = "22.412.21"
$ws.Range("E25").Value = "  -4.64%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "2.403"
$ws.Range("E26").Value = "  -4.45%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "2.949"
$ws.Range("E27").Value = "  +1.11%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "19.81"
$ws.Range("E28").Value = "  -5.24%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "146.29"
$ws.Range("E29").Value = "  -4.69%  "
$ws.Range("B30").Value = "HuobiToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D30").Value = "4.921"
$ws.Range("E30").Value = "  -5.78%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "125.01"
$ws.Range("E31").Value = "  -5.90%  "
$ws.Range("B32").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C32").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D32").Value = "1.751.64"
$ws.Range("E32").Value = "  -4.81%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "6.256"
$ws.Range("E33").Value = "  -9.77%  "
$ws.Range("B34").Value = "WEMIXTOKEN"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Value = "1.981"
$ws.Range("E34").Value = "  -5.88%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "0.9822"
$ws.Range("E35").Value = "  -3.42%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "10.35"
$ws.Range("E36").Value = "  -13.39%  "
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").Value = "0.08425"
$ws.Range("E37").Value = "  -3.56%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.02531"
$ws.Range("E38").Value = "  -7.30%  "
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "0.2298"
$ws.Range("E39").Value = "  -5.74%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.06504"
$ws.Range("E40").Value = "  -4.36%  "
$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").Value = "5.492"
$ws.Range("E41").Value = "  -7.73%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "11.77"
$ws.Range("E42").Value = "  -10.53%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "1.242"
$ws.Range("E43").Value = "  -5.64%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.6387"
$ws.Range("E44").Value = "  -7.52%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "14.55"
$ws.Range("E45").Value = "  -6.22%  "
$ws.Range("B46").Value = "Frax"
$ws.Range("C46").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.6008"
$ws.Range("E47").Value = "  -6.43%  "
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "3.773"
$ws.Range("E48").Value = "  -3.81%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "2.106"
$ws.Range("E49").Value = "  -7.00%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "122.12"
$ws.Range("E50").Value = "  -4.32%  "
$ws.Range("D51").Value = "1.190"
$ws.Range("E51").Value = "  +0.52%  "
